# Update the "想去人数" (attendance/interest count) figures in column F
# for the worksheets that list event data: "展览" and "全部类型".
# Both sheets mirror the same dataset, so the same cell updates apply to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F3"  = 7180
    "F4"  = 5272
    "F6"  = 165
    "F9"  = 104
    "F10" = 76
    "F11" = 93
    "F13" = 635
    "F14" = 225
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
